# Actualización automática 2025-06-30 12:55:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("F4").Value = 76.75
$wsVentasGrupo.Range("E16").Value = 401.84
$wsVentasGrupo.Range("E55").Value = "5 de 53"
$wsVentasGrupo.Range("F55").Value = "1 de 53"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F4").Value = 24916.2
$wsVentaMensual.Range("F16").Value = 5554.89
$wsVentaMensual.Range("F55").Value = 77844.75

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D4").Value = 1568.16
$wsCumplimiento.Range("E4").Value = -565.1600000000001
$wsCumplimiento.Range("F4").Value = 1.563469591226321

$wsCumplimiento.Range("D5").Value = 76.75
$wsCumplimiento.Range("E5").Value = 161.57
$wsCumplimiento.Range("F5").Value = 0.3220459885867741

$wsCumplimiento.Range("D19").Value = 80259.67999999999
$wsCumplimiento.Range("E19").Value = 14187.76064517915
$wsCumplimiento.Range("F19").Value = 0.8497814175983885
